$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.956.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.88%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.663.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.29%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'523.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.19%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'144.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.00%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.20%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.570"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.40%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'6.98"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +7.35%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -3.64%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -2.16%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D13").Value = "'3.133.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.78%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'58.946.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.99%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -1.94%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'ShibaInu"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.0000136"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.25%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'WrappedEther"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'2.657.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -5.89%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'338.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.18%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -3.74%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -2.78%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.11%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.08%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'64.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.46%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -1.22%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -1.44%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.27%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -2.84%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.28%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'6.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.09%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.07%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.45%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'18.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.36%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'150.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.71%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.11%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -5.57%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.899"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -6.22%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.870"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'36.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.41%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -5.93%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -3.58%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.616"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.10%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.11%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'275.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.69%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'19.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.69%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0969"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.33%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'10.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.94%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -1.43%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.051.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.55%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'4.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.17%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -3.15%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'18.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.25%  "
$ws.Range("E51").Style = "Normal"
